# Adding test Case to Search Module OPQA_1243
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new shared-string values in the exact order they need to be
# interned so the shared-strings table matches the target layout
# (B114, C114, B115, C115, A114, A115 introduce new strings; the D/E
# values reuse existing "Y" / "PASS" strings already in the table).
$ws.Range("B114").Value = "OPQA-593"
$ws.Range("C114").Value = "Verify that MORE and LESS links are working correctly in INVENTOR filter in PATENTS search results page"
$ws.Range("B115").Value = "OPQA-588"
$ws.Range("C115").Value = "Verify that left navigation pane content type is retained when user navigates back to PATENTS search results page from record view page"
$ws.Range("A114").Value = "TestCase_B113"
$ws.Range("A115").Value = "TestCase_B114"
$ws.Range("D114").Value = "Y"
$ws.Range("E114").Value = "PASS"
$ws.Range("D115").Value = "Y"
$ws.Range("E115").Value = "PASS"

# Copy the formatting of the last existing data row (113) onto the two
# new rows so they reuse the same thin-border cell style rather than a
# brand new style definition being created.
$ws.Range("A113:E113").Copy()
$ws.Range("A114:E114").PasteSpecial(-4122)
$ws.Range("A113:E113").Copy()
$ws.Range("A115:E115").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to reflect the new extent of data (full table
# selected). The engine's Range.Activate() collapses a multi-cell
# selection down to the single activated cell, so select the full range
# last to keep sqref spanning A1:E115 (matching the new data extent).
$ws.Range("A1:E115").Select()
